$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Product_Name/chr/Name of Pizza-Beverage -> margherita_pizza_price/num/Price of 32cm Margherita pizza in HUF
$ws.Range("A3").Value = "margherita_pizza_price"
$ws.Range("B3").Value = "num"
$ws.Range("C3").Value = "Price of 32cm Margherita pizza in HUF"

# Row 4: Price/num/Price in HUF -> beverage_price/num/Price of 0.5 liter beverage in HUF (only Pepsi and Coca-Cola)
$ws.Range("A4").Value = "beverage_price"
$ws.Range("B4").Value = "num"
$ws.Range("C4").Value = "Price of 0.5 liter beverage in HUF (only Pepsi and Coca-Cola)"

# Row 16: Size/num/Size of Pizza in cm or Beverage in l -> Tags/chr/Concatenated string...
$ws.Range("A16").Value = "Tags"
$ws.Range("B16").Value = "chr"
$ws.Range("C16").Value = "Concatenated string with max 5 tags available for a restaurant "

# Delete the last row (previously row 17, Tags duplicate) and shift cells up.
$ws.Range("A17:C17").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

$wb.Save()
